$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "35÷3=" -> "25÷2="
$cell = $t.Cell(1, 1)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "25÷2="

# Row 1, Col 2: "75÷2=" -> "10÷4="
$cell = $t.Cell(1, 2)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "10÷4="

# Row 1, Col 3: "72÷4=" -> "42÷8="
$cell = $t.Cell(1, 3)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "42÷8="

# Row 1, Col 4: "51÷4=" -> "93÷7="
$cell = $t.Cell(1, 4)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "93÷7="

# Row 1, Col 5: "83÷2=" -> "39÷3="
$cell = $t.Cell(1, 5)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "39÷3="

# Row 5, Col 1: "90÷7=" -> "98÷7="
$cell = $t.Cell(5, 1)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "98÷7="

# Row 5, Col 2: "48÷4=" -> "58÷8="
$cell = $t.Cell(5, 2)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "58÷8="

# Row 5, Col 3: "89÷8=" -> "15÷3="
$cell = $t.Cell(5, 3)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "15÷3="

# Row 5, Col 4: "40÷7=" -> "95÷8="
$cell = $t.Cell(5, 4)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "95÷8="

# Row 5, Col 5: "31÷6=" -> "48÷2="
$cell = $t.Cell(5, 5)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "48÷2="

# Row 9, Col 1: "15÷5=" -> "27÷8="
$cell = $t.Cell(9, 1)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "27÷8="

# Row 9, Col 2: "90÷8=" -> "63÷3="
$cell = $t.Cell(9, 2)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "63÷3="

# Row 9, Col 3: "44÷7=" -> "69÷2="
$cell = $t.Cell(9, 3)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "69÷2="

# Row 9, Col 4: "48÷4=" -> "70÷6="
$cell = $t.Cell(9, 4)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "70÷6="

# Row 9, Col 5: "76÷4=" -> "91÷9="
$cell = $t.Cell(9, 5)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "91÷9="

# Row 13, Col 1: "24÷4=" -> "15÷2="
$cell = $t.Cell(13, 1)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "15÷2="

# Row 13, Col 2: "44÷3=" -> "21÷5="
$cell = $t.Cell(13, 2)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "21÷5="

# Row 13, Col 3: "88÷2=" -> "56÷2="
$cell = $t.Cell(13, 3)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "56÷2="

# Row 13, Col 4: "93÷3=" -> "93÷2="
$cell = $t.Cell(13, 4)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "93÷2="

# Row 13, Col 5: "45÷6=" -> "73÷4="
$cell = $t.Cell(13, 5)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "73÷4="

# Row 17, Col 1: "39÷5=" -> "34÷9="
$cell = $t.Cell(17, 1)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "34÷9="

# Row 17, Col 2: "72÷4=" -> "34÷6="
$cell = $t.Cell(17, 2)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "34÷6="

# Row 17, Col 3: "53÷7=" -> "18÷6="
$cell = $t.Cell(17, 3)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "18÷6="

# Row 17, Col 4: "55÷2=" -> "64÷9="
$cell = $t.Cell(17, 4)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "64÷9="

# Row 17, Col 5: "75÷3=" -> "36÷6="
$cell = $t.Cell(17, 5)
$r = $cell.Range
$r.SetRange($r.Start, $r.End - 1)
$r.Text = "36÷6="
